$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

$ws.Range("C3").Value = "Done"
$ws.Range("C8").Value = "Done"
$ws.Range("C9").Value = "Done"

$ws.Range("C9").Select()
